$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top and set its value, mirroring a user
# right-clicking the row-1 header and choosing Insert, then typing
# a header label into A1.
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "Adresy"

# Leave the whole of row 1 selected, like after a row-header insert.
$ws.Rows.Item(1).Select() | Out-Null
